$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 585, shifting all following
# rows (585..621) down to (586..622). Excel copies the formatting of the
# row above into the new row, matching the original D-column date style.
$ws.Rows.Item(585).Insert()

# Populate the newly-inserted row 585 with the new weekly record.
$ws.Range("A585").Value = 7
$ws.Range("B585").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C585").Value = "Ñuble"
$ws.Range("D585").Value = 44578
$ws.Range("E585").Value = 16
$ws.Range("F585").Value = "Fruta"
$ws.Range("G585").Value = 100104
$ws.Range("H585").Value = "Frutos de pepita"
$ws.Range("I585").Value = 100104002
$ws.Range("J585").Value = "Manzana"
$ws.Range("K585").Value = "Fuji royal"
$ws.Range("L585").Value = "Primera"
$ws.Range("M585").Value = 100
$ws.Range("N585").Value = 12000
$ws.Range("O585").Value = 13000
$ws.Range("P585").Value = 12500
$ws.Range("Q585").Value = "$/caja 16 kilos empedrada"
$ws.Range("R585").Value = "Provincia de Curicó"
$ws.Range("S585").Value = 781
$ws.Range("T585").Value = 16
